$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the weight matrix values
$ws.Range("D2").Value = 0.11
$ws.Range("D3").Value = 0.3
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 3

# Update the selected cell to match the saved selection state
$ws.Range("B6").Select()
